$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new wishlist entry as row 9
$ws.Range("A9").Value = "La ciudad de las bestias"
$ws.Range("B9").Value = "Isabel Allende"

# C8 is an existing empty cell (present but with no value). Copy it down to
# C9 so the new row also gets an empty-but-present cell in column C,
# matching the pattern of the other rows.
$ws.Range("C8").Copy($ws.Range("C9"))
